$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 9 (item 6) — new non-conformity about physical location of files not
# being respected (recent projects only)
# ---------------------------------------------------------------------------
$ws.Range("B8").Copy()
$ws.Range("B9").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B9").Value = "No se esta respetando las ubicación física de los archivos establecido por el plan de configuración (solo se presento en proyectos recientes)"

$ws.Range("C9").Value = "Equipo de ventas"

$ws.Range("D5").Copy()
$ws.Range("D9").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D9").Value = 42460

$ws.Range("F9").Value = "En proceso"

$ws.Rows(9).RowHeight = 39.55

# ---------------------------------------------------------------------------
# Row 10 (item 7) — satisfaction surveys not carried out
# ---------------------------------------------------------------------------
$ws.Range("B10").Value = "No se tiene realizadas encuestas de satisfacción"
$ws.Range("C10").Value = "Magda Montoya"

$ws.Range("D5").Copy()
$ws.Range("D10").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D10").Value = 42460

$ws.Range("F10").Value = "En proceso"

$ws.Rows(10).RowHeight = 14.9

# ---------------------------------------------------------------------------
# Row 11 (item 8) — new non-conformity about support tickets lacking time
# tracking
# ---------------------------------------------------------------------------
$ws.Range("B11").Value = "Los tickets de soporte no cuentan con un tracking de tiempo"
$ws.Range("C11").Value = "Equipo de ventas"

$ws.Range("D5").Copy()
$ws.Range("D11").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D11").Value = 42460

$ws.Range("F11").Value = "En proceso"

$ws.Rows(11).RowHeight = 28.35

# ---------------------------------------------------------------------------
# Row 12 (item 9) — still empty, only the row height shrinks back down
# ---------------------------------------------------------------------------
$ws.Rows(12).RowHeight = 13.8

# ---------------------------------------------------------------------------
# Selection moved from F8 to B10, with the view scrolled back up to row 4
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("B10").Select()
